$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last (duplicate) data row - row 5 was an exact copy of row 4.
$ws.Rows(5).Delete()

# Insert a new blank column before column E; this shifts the old E:L block
# (the SE statistic columns) one column to the right, into F:M.
$ws.Columns(5).Insert()

# Row 4 gains a label in the newly inserted column E.
$ws.Range("E4").Value = "Var"

# Refresh the numeric results in the shifted SE columns (F:M) to the
# recomputed values.
$row2 = @(0.47, 0.36, 1, 0.08, 0.2, 0.5, 0.84, 0.25)
$row3 = @(0.47, 0.38, 1, 0.1, 0.21, 0.54, 0.92, 0.18)
$row4 = @(0.47, 0.36, 1, 0.08, 0.21, 0.5, 0.84, 0.25)
$cols = @("F", "G", "H", "I", "J", "K", "L", "M")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $row2[$i]
    $ws.Range($cols[$i] + "3").Value = $row3[$i]
    $ws.Range($cols[$i] + "4").Value = $row4[$i]
}

# Leave the selection on the last updated cell, matching the author's
# final cursor position.
$ws.Range("M4").Select() | Out-Null
